# Fix repeated demo number in presentation.
#
# Slides 12 and 13 (the "Setting up / configuring BioJava" slides under
# PART 1) both incorrectly carried the title "Demo 2" -- duplicating the
# title that legitimately belongs to the later structures demo (slide 18,
# under PART 2). Retitle slides 12 and 13 to "Setting up BioJava",
# matching the existing "Setting up BioJava" title used on slide 6, and
# split the run the same way ("Setting up " + "BioJava") since BioJava
# is flagged by the spell checker elsewhere in the deck.

$p = $ppt.ActivePresentation

foreach ($idx in 12, 13) {
    $slide = $p.Slides.Item($idx)
    $title = $slide.Shapes.Item(1).TextFrame.TextRange

    $title.Text = "Setting up BioJava"

    # Split into two runs ("Setting up " / "BioJava"), mirroring slide 6's
    # title so BioJava keeps its own run (it gets flagged by the spell
    # checker as a foreign/unknown word there too).
    $bioJava = $title.Characters(12, 8)
    $bioJava.Text = "BioJava"
}
